$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 264.375
$ws.Range("I2").Value = 273.57144
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 273.57144
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = -160.57144
$ws.Range("N2").Value = -426
$ws.Range("H9").Value = 148.44444
$ws.Range("J9").Value = 99
$ws.Range("L9").Value = 99
$ws.Range("N9").Value = -437
$ws.Range("H29").Value = 3975.1365
$ws.Range("I29").Value = 3194.0557
$ws.Range("J29").Value = 7490
$ws.Range("K29").Value = 9582.167099999999
$ws.Range("L29").Value = 22470
$ws.Range("M29").Value = -9301.167099999999
$ws.Range("N29").Value = -23032
$ws.Range("H38").Value = 3661.4375
$ws.Range("I38").Value = 2068.5557
$ws.Range("J38").Value = 5709.4287
$ws.Range("K38").Value = 6205.6671
$ws.Range("L38").Value = 17128.2861
$ws.Range("M38").Value = -5833.6671
$ws.Range("N38").Value = -17872.2861
$ws.Range("H43").Value = 8875
$ws.Range("J43").Value = 11166.667
$ws.Range("L43").Value = 11166.667
$ws.Range("N43").Value = -11304.667
$ws.Range("H53").Value = 355.66666
$ws.Range("J53").Value = 178.75
$ws.Range("L53").Value = 178.75
$ws.Range("N53").Value = -1452.75
$ws.Range("H58").Value = 843.4286
$ws.Range("I58").Value = 150.66667
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 452.00001
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = -302.00001
$ws.Range("N58").Value = -15300
$ws.Range("H116").Value = 2215.7693
$ws.Range("I116").Value = 2079.4
$ws.Range("K116").Value = 2079.4
$ws.Range("M116").Value = 1362.6
$ws.Range("H137").Value = 1854.5
$ws.Range("I137").Value = 1225.4
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 3676.2
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = -1126.2
$ws.Range("N137").Value = -20100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 1196542.6
$ws.Range("J135").Value = 1196542.6
$ws.Range("L135").Value = 1196542.6
$ws.Range("N135").Value = -1206682.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4354.0835
$ws.Range("I58").Value = 3615.5715
$ws.Range("K58").Value = 3615.5715
$ws.Range("M58").Value = -3412.5715
$ws.Range("H62").Value = 5249.5
$ws.Range("J62").Value = 1998
$ws.Range("L62").Value = 1998
$ws.Range("N62").Value = -3246
$ws.Range("H65").Value = 5249.5
$ws.Range("J65").Value = 1998
$ws.Range("L65").Value = 9990
$ws.Range("N65").Value = -16230
$ws.Range("H122").Value = 2362.5
$ws.Range("I122").Value = 2362.5
$ws.Range("K122").Value = 7087.5
$ws.Range("M122").Value = -4637.5
$ws.Range("H136").Value = 4354.0835
$ws.Range("I136").Value = 3615.5715
$ws.Range("K136").Value = 10846.7145
$ws.Range("M136").Value = -8296.7145

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3899.9
$ws.Range("J139").Value = 6000
$ws.Range("L139").Value = 18000
$ws.Range("N139").Value = -28280

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2795.75
$ws.Range("I80").Value = 2647
$ws.Range("J80").Value = 2944.5
$ws.Range("K80").Value = 2647
$ws.Range("L80").Value = 2944.5
$ws.Range("M80").Value = -1649
$ws.Range("N80").Value = -4940.5
$ws.Range("H83").Value = 2795.75
$ws.Range("I83").Value = 2647
$ws.Range("J83").Value = 2944.5
$ws.Range("K83").Value = 13235
$ws.Range("L83").Value = 14722.5
$ws.Range("M83").Value = -8243
$ws.Range("N83").Value = -24706.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6079.0835
$ws.Range("I7").Value = 5822.25
$ws.Range("J7").Value = 6592.75
$ws.Range("K7").Value = 5822.25
$ws.Range("L7").Value = 6592.75
$ws.Range("M7").Value = -5710.25
$ws.Range("N7").Value = -6816.75
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -705
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -893
$ws.Range("N27").ClearContents()
$ws.Range("H46").Value = 3907.3333
$ws.Range("I46").Value = 1696
$ws.Range("J46").Value = 8330
$ws.Range("K46").Value = 1696
$ws.Range("L46").Value = 8330
$ws.Range("M46").Value = -1508
$ws.Range("N46").Value = -8706
$ws.Range("H68").Value = 10000
$ws.Range("J68").Value = 10000
$ws.Range("L68").Value = 10000
$ws.Range("N68").Value = -11498
$ws.Range("H71").Value = 10000
$ws.Range("J71").Value = 10000
$ws.Range("L71").Value = 50000
$ws.Range("N71").Value = -57488
$ws.Range("H82").Value = 5149.8887
$ws.Range("I82").Value = 3000
$ws.Range("J82").Value = 5418.625
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 5418.625
$ws.Range("M82").Value = -2639
$ws.Range("N82").Value = -6140.625
$ws.Range("H85").Value = 5149.8887
$ws.Range("I85").Value = 3000
$ws.Range("J85").Value = 5418.625
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 5418.625
$ws.Range("M85").Value = -1752
$ws.Range("N85").Value = -7914.625
$ws.Range("H126").Value = 6079.0835
$ws.Range("I126").Value = 5822.25
$ws.Range("J126").Value = 6592.75
$ws.Range("K126").Value = 17466.75
$ws.Range("L126").Value = 19778.25
$ws.Range("M126").Value = -14996.75
$ws.Range("N126").Value = -24718.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11500
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 11500
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 11500
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -12748
$ws.Range("H65").Value = 11500
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 11500
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 57500
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -63740
$ws.Range("J126").Value = 7690
$ws.Range("L126").Value = 23070
$ws.Range("N126").Value = -28010
